# Insert a new timestamped snapshot column ("EI") in the price-history
# sheet, shifting the trailing "nom" / "url_produit" columns one slot to
# the right (EI->EJ, EJ->EK), and populate the new column:
#   - row 1 (header): the new scrape timestamp
#   - rows 2..206: the same value that was already in column EH (the
#     previous snapshot) for that row - i.e. carry the last known price
#     forward (or leave blank for rows that had no price yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 206
$oldLastCol = 138   # EH: last snapshot column before the insert

# Shift EI:EJ -> EJ:EK, leaving a fresh empty column at EI.
$ws.Columns("EI").Insert()

# New header timestamp for the inserted column.
$ws.Cells.Item(1, 139).Value = "2026-02-03 08:24:36"

# Carry forward column EH's value into the newly inserted EI column for
# every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $prev = $ws.Cells.Item($r, $oldLastCol).Value2
    if ($prev -ne "") {
        $ws.Cells.Item($r, 139).Value = $prev
    }
}

Write-Host "done"
